$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.948.92"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "3.690.32"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "652.79"
$ws.Range("E5").Value = "  -3.63%  "

$ws.Range("D6").Value = "161.99"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.499"
$ws.Range("E8").Value = "  +0.58%  "

$ws.Range("E9").Value = "  -1.51%  "

$ws.Range("D10").Value = "7.16"
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("D11").Value = "0.443"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").Value = "4.314.02"
$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").Value = "32.76"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").Value = "3.710.07"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").Value = "69.949.36"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").Value = "16.01"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "6.54"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").Value = "10.39"
$ws.Range("E20").Value = "  +6.00%  "

$ws.Range("D21").Value = "471.61"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").Value = "79.87"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").Value = "3.840.07"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("E25").Value = "  +0.39%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").Value = "11.06"
$ws.Range("E27").Value = "  +0.94%  "

$ws.Range("D28").Value = "8.87"
$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("E29").Value = "  -1.74%  "

$ws.Range("D30").Value = "1.71"
$ws.Range("E30").Value = "  -2.13%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("D33").Value = "6.54"
$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").Value = "26.78"
$ws.Range("E35").Value = "  -0.68%  "

$ws.Range("D36").Value = "3.686.09"
$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("D37").Value = "8.42"
$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "5.91"
$ws.Range("E39").Value = "  -4.95%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "178.49"
$ws.Range("E40").Value = "  +7.04%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -1.09%  "

$ws.Range("D43").Value = "0.0902"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "0.931"
$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("D45").Value = "2.84"
$ws.Range("E45").Value = "  +2.30%  "

$ws.Range("D46").Value = "29.28"
$ws.Range("E46").Value = "  +3.13%  "

$ws.Range("D47").Value = "46.58"

$ws.Range("D48").Value = "0.000273"
$ws.Range("E48").Value = "  -2.69%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "7.88"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.25"
$ws.Range("E50").Value = "  -3.48%  "

$ws.Range("D51").Value = "1.05"
$ws.Range("E51").Value = "  -4.99%  "
